# orm-03.pptx edit script
# 1) Remove the extra "Demo" slide (old slide #6 / SlideID 300). Its notes
#    slide (notesSlide2.xml) and associated relationships are cleaned up
#    automatically by the host when the slide is deleted.
$p = $ppt.ActivePresentation

$p.Slides.Item(6).Delete()

# 2) The slide that used to be #7 ("Exercise") is now #6 and its cached
#    slide-number field needs to read 6 instead of 7. Toggling the
#    Slide Number header/footer visibility off/on forces the host to
#    recompute the cached field text for the slide's current position.
$sExercise = $p.Slides.Item(6)
$sExercise.HeadersFooters.SlideNumber.Visible = $false
$sExercise.HeadersFooters.SlideNumber.Visible = $true

# 3) Slide 4 ("Demo" / find options) - switch the "select"/"relations"/"join"
#    bullets (level 3) to the Consolas monospace font.
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
foreach ($i in 5..7) {
    $tr4.Paragraphs($i, 1).Font.Name = "Consolas"
}

# 4) Slide 5 ("Demo" / find options continued) - apply Consolas to the
#    where/and/or/order/skip,take/Not,LessThan,... bullets, bump "and"/"or"
#    to 18pt, and fix up a few labels (TypeORM renamed these operators).
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

$tr5.Paragraphs(3, 1).Font.Name = "Consolas"                 # where

$pAnd = $tr5.Paragraphs(4, 1)                                 # and
$pAnd.Font.Name = "Consolas"
$pAnd.Font.Size = 18

$pOr = $tr5.Paragraphs(5, 1)                                  # or
$pOr.Font.Name = "Consolas"
$pOr.Font.Size = 18

$tr5.Paragraphs(6, 1).Font.Name = "Consolas"                  # order

$pSkip = $tr5.Paragraphs(7, 1)                                # skip, limit -> skip, take
$pSkip.Characters(1, $pSkip.Length).Text = "skip, take"
$tr5.Paragraphs(7, 1).Font.Name = "Consolas"

$pLast = $tr5.Paragraphs(8, 1)
# "Not, " "lessThan" ", " "moreThan" ", equal, Like, between, In, Raw"
$pLast.Characters(1, 5).Font.Name = "Consolas"                # "Not, "

$rLess = $pLast.Characters(6, 8)                              # "lessThan" -> "LessThan"
$rLess.Text = "LessThan"
$rLess.Font.Name = "Consolas"

$pLast.Characters(14, 2).Font.Name = "Consolas"               # ", "

$rMore = $pLast.Characters(16, 8)                             # "moreThan" -> "MoreThan"
$rMore.Text = "MoreThan"
$rMore.Font.Name = "Consolas"

$rTail = $pLast.Characters(24, 31)                            # ", equal, ..." -> ", Equal, ..."
$rTail.Text = ", Equal, Like, between, In, Raw"
$rTail.Font.Name = "Consolas"
